# Update cryptos list with latest prices/volumes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to stay text even when the string looks like a pure
# number (e.g. "250.07"), mirroring how the source data is stored as text.
function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Column D price updates that are unambiguous numeric strings ---
# (these need the text NumberFormat so Excel doesn't coerce them to Double)
Set-TextCell "D5"  "250.07"
Set-TextCell "D7"  "61.61"
Set-TextCell "D10" "0.0794"
Set-TextCell "D12" "16.33"
Set-TextCell "D14" "0.829"
Set-TextCell "D17" "18.34"
Set-TextCell "D19" "75.34"
Set-TextCell "D21" "5.45"
Set-TextCell "D22" "239.46"
Set-TextCell "D25" "2.22"
Set-TextCell "D26" "169.51"
Set-TextCell "D28" "20.06"
Set-TextCell "D34" "0.0892"
Set-TextCell "D35" "1.00"
Set-TextCell "D42" "18.20"
Set-TextCell "D43" "0.0225"
Set-TextCell "D44" "1.15"
Set-TextCell "D45" "98.18"
Set-TextCell "D46" "2.50"
Set-TextCell "D48" "2.89"
Set-TextCell "D49" "6.87"
Set-TextCell "D51" "3.54"

# --- Remaining plain text updates (prices with thousand separators, all
#     Volume(1h) percentages, and the two swapped coin names/links) ---

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.171.84"
$ws.Range("E2").Value = "  +0.15%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.057.31"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.34%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.22%  "

# Row 7 - Solana
$ws.Range("E7").Value = "  +10.22%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.11%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.53%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.92%  "

# Row 12 - Chainlink
$ws.Range("E12").Value = "  +7.95%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.355.72"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +1.36%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +8.90%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.055.02"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17 - Avalanche
$ws.Range("E17").Value = "  +28.66%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.181.99"
$ws.Range("E18").Value = "  +0.37%  "

# Row 19 - Litecoin
$ws.Range("E19").Value = "  +3.80%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -4.83%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.61%  "

# Row 22 - BitcoinCash
$ws.Range("E22").Value = "  +0.86%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.00%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.95%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +12.18%  "

# Row 26 - Monero
$ws.Range("E26").Value = "  -0.69%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +4.13%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -0.40%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +1.75%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +9.97%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +5.70%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -0.31%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +4.19%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  +4.60%  "

# Row 35 - BinanceUSD
$ws.Range("E35").Value = "  -0.09%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -0.91%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  -2.23%  "

# Row 38 - Cronos
$ws.Range("E38").Value = "  +2.16%  "

# Row 39 - TrustWalletToken
$ws.Range("E39").Value = "  +0.61%  "

# Row 40 - THORChain
$ws.Range("E40").Value = "  +31.35%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +14.73%  "

# Row 42 - InjectiveProtocol
$ws.Range("E42").Value = "  +0.67%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +0.40%  "

# Row 44 & 45 - Aave and ARBITRUM swap positions (with refreshed data)
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E45").Value = "  +1.62%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +3.20%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.298.43"

# Row 48 - MXToken
$ws.Range("E48").Value = "  -1.10%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  +0.47%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.242.46"
$ws.Range("E50").Value = "  -0.39%  "

# Row 51 - FTXToken
$ws.Range("E51").Value = "  -16.50%  "
